# ---------------------------------------------------------------------------
# Insert a new "2022-Q3" sheet (fund detail breakdown) right after "总计",
# and update the "总计" (totals) sheet with the new quarter's summary row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$wsTotal = $wb.Worksheets.Item(1)          # 总计

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsTotal)
$wsQ3.Name = "2022-Q3"
$wsQ3.Activate()

# Fetch the (still-existing) "2022-Q2" sheet AFTER the insertion above, since
# worksheet references taken before a Worksheets.Add() call can go stale.
$wsQ2 = $wb.Worksheets.Item("2022-Q2")     # existing sheet used as a style template

# Copy the header-row / column-A cell styles from the existing "2022-Q2" sheet
# so the new sheet matches the look (bold header + bordered column A).
$wsQ2.Range("B1:H1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$wsQ2.Range("A2").Copy()
$wsQ3.Range("A2:A10").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------------
# Header row
# ---------------------------------------------------------------------------
$wsQ3.Range("B1").Value = "基金代码"
$wsQ3.Range("C1").Value = "基金名称"
$wsQ3.Range("D1").Value = "基金规模"
$wsQ3.Range("E1").Value = "股票总仓位"
$wsQ3.Range("F1").Value = "仓位占比"
$wsQ3.Range("G1").Value = "持有市值(亿元)"
$wsQ3.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# Data rows (A = index, B..G stored as text to mirror the source data,
# H stored as a real number)
# ---------------------------------------------------------------------------
$rows = @(
    @{ A=0; B="016935"; C="景顺长城中证500指数增强C";                           D="15.57"; E="93.89"; F="1.75"; G="0.2725"; H=9  },
    @{ A=1; B="000978"; C="景顺长城量化精选股票";                               D="7.14";  E="93.64"; F="1.84"; G="0.1314"; H=9  },
    @{ A=2; B="160218"; C="国泰国证房地产行业指数A";                           D="5.91";  E="94.23"; F="1.93"; G="0.1141"; H=10 },
    @{ A=3; B="515760"; C="华夏中证浙江国资创新发展ETF";                       D="2.04";  E="99.57"; F="3.39"; G="0.0692"; H=9  },
    @{ A=4; B="008851"; C="景顺长城量化对冲策略三个月定期开放灵活配置混合";     D="2.96";  E="64.77"; F="1.34"; G="0.0397"; H=6  },
    @{ A=5; B="015042"; C="国泰国证房地产行业指数C";                           D="1.29";  E="94.23"; F="1.93"; G="0.0249"; H=10 },
    @{ A=6; B="004157"; C="信诚至诚灵活配置混合A";                             D="0.63";  E="24.85"; F="0.76"; G="0.0048"; H=5  },
    @{ A=7; B="004158"; C="信诚至诚灵活配置混合B";                             D="0.17";  E="24.85"; F="0.76"; G="0.0013"; H=5  },
    @{ A=8; B="006682"; C="景顺长城中证500指数增强A";                          D="0.00";  E="93.89"; F="1.75"; G=0;        H=9  }
)

$r = 2
foreach ($row in $rows) {
    $wsQ3.Range("A$r").Value = $row.A

    $wsQ3.Range("B$r").Value = "'" + $row.B
    $wsQ3.Range("B$r").Style = "Normal"

    $wsQ3.Range("C$r").Value = "'" + $row.C
    $wsQ3.Range("C$r").Style = "Normal"

    $wsQ3.Range("D$r").Value = "'" + $row.D
    $wsQ3.Range("D$r").Style = "Normal"

    $wsQ3.Range("E$r").Value = "'" + $row.E
    $wsQ3.Range("E$r").Style = "Normal"

    $wsQ3.Range("F$r").Value = "'" + $row.F
    $wsQ3.Range("F$r").Style = "Normal"

    if ($row.G -eq 0) {
        $wsQ3.Range("G$r").Value = 0
    } else {
        $wsQ3.Range("G$r").Value = "'" + $row.G
        $wsQ3.Range("G$r").Style = "Normal"
    }

    $wsQ3.Range("H$r").Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2. Update the "总计" (totals) sheet: insert a new top row for 2022-Q3 and
#    push the other quarters' rows down by one.
# ---------------------------------------------------------------------------
$wsTotal.Range("A5").Copy()
$wsTotal.Range("A6").PasteSpecial(-4122)   # xlPasteFormats
$wsTotal.Range("A6").Value = 4
$wsTotal.Range("B6").Value = "2021-Q2"
$wsTotal.Range("C6").Value = 1
$wsTotal.Range("D6").Value = 0.13

$wsTotal.Range("B5").Value = "2021-Q4"
$wsTotal.Range("C5").Value = 1
$wsTotal.Range("D5").Value = 0.08

$wsTotal.Range("B4").Value = "2022-Q1"
$wsTotal.Range("C4").Value = 2
$wsTotal.Range("D4").Value = 0.13

$wsTotal.Range("B3").Value = "2022-Q2"
$wsTotal.Range("C3").Value = 9
$wsTotal.Range("D3").Value = 0.42

$wsTotal.Range("B2").Value = "2022-Q3"
$wsTotal.Range("C2").Value = 9
$wsTotal.Range("D2").Value = 0.66
